{"js": "// Update the recalibrated performance-metric codes in the\n// \"Imputed Performance metB\" table (commit: \"Update codes for recalibration\").\n//\n// Each (row, column) pair below addresses a specific table cell in the\n// single table that exists in this document. Column 2 (\"Method B\") holds\n// the numeric/interval values that changed. Addressing cells by their\n// table position (rather than a blind, document-wide text search-and-\n// replace) avoids any ambiguity from values that coincidentally collide\n// with each other after the edit (e.g. the new \"Calibration Intercept\"\n// value at 5 years equals the *old* \"Calibration Slope\" value at 2 years).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document, found none.\");\n}\n\nconst table = tables.items[0];\n\n// [rowIndex, columnIndex, oldText, newText]\nconst edits = [\n  [4, 2, \"2.73%\", \"2.82%\"],\n  [5, 2, \"2.6% (2.41% to 2.78%)\", \"2.73% (2.54% to 2.92%)\"],\n  [6, 2, \"0.95 (0.88 to 1.02)\", \"0.97 (0.9 to 1.04)\"],\n  [7, 2, \"-0.14% (-0.32% to 0.05%)\", \"-0.09% (-0.28% to 0.1%)\"],\n  [8, 2, \"-0.06 (-0.16 to 0.04)\", \"-0.05 (-0.16 to 0.05)\"],\n  [9, 2, \"-0.14 (-0.21 to -0.07)\", \"-0.15 (-0.24 to -0.06)\"],\n  [11, 2, \"0.91 (0.89 to 0.92)\", \"0.9 (0.89 to 0.92)\"],\n  [17, 2, \"5.02%\", \"5.23%\"],\n  [18, 2, \"4.49% (4.24% to 4.74%)\", \"4.76% (4.51% to 5.02%)\"],\n  [19, 2, \"0.89 (0.84 to 0.94)\", \"0.91 (0.86 to 0.96)\"],\n  [20, 2, \"-0.54% (-0.79% to -0.28%)\", \"-0.47% (-0.72% to -0.22%)\"],\n  [21, 2, \"-0.18 (-0.26 to -0.1)\", \"-0.14 (-0.21 to -0.07)\"],\n  [22, 2, \"-0.16 (-0.22 to -0.1)\", \"-0.16 (-0.21 to -0.1)\"],\n  [24, 2, \"0.89 (0.88 to 0.9)\", \"0.88 (0.87 to 0.89)\"],\n  [26, 2, \"0.04 (0.03 to 0.04)\", \"0.04 (0.04 to 0.04)\"],\n];\n\nfor (const [rowIndex, colIndex, oldText, newText] of edits) {\n  const cell = table.getCell(rowIndex, colIndex);\n  // Scope the search to this single cell's body so the replacement can\n  // never touch a different cell, even if two values happen to coincide\n  // after earlier edits in this loop run.\n  const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\n      `Could not find expected text \"${oldText}\" in table cell (row ${rowIndex}, col ${colIndex}).`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the recalibrated performance-metric codes in the\n# \"Imputed Performance metB\" table (commit: \"Update codes for recalibration\").\n#\n# Each entry addresses a specific table cell by its 1-based (row, column)\n# position in the document's single table. Column 3 (\"Method B\") holds the\n# numeric/interval values that changed. Addressing cells by their table\n# position (rather than a blind, document-wide Find/Replace) avoids any\n# ambiguity from values that coincidentally collide with each other after\n# the edit (e.g. the new \"Calibration Intercept\" value at 5 years equals\n# the *old* \"Calibration Slope\" value at 2 years).\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# row, col (1-based), expected old text, new text\n$edits = @(\n    @(5, 3, \"2.73%\", \"2.82%\"),\n    @(6, 3, \"2.6% (2.41% to 2.78%)\", \"2.73% (2.54% to 2.92%)\"),\n    @(7, 3, \"0.95 (0.88 to 1.02)\", \"0.97 (0.9 to 1.04)\"),\n    @(8, 3, \"-0.14% (-0.32% to 0.05%)\", \"-0.09% (-0.28% to 0.1%)\"),\n    @(9, 3, \"-0.06 (-0.16 to 0.04)\", \"-0.05 (-0.16 to 0.05)\"),\n    @(10, 3, \"-0.14 (-0.21 to -0.07)\", \"-0.15 (-0.24 to -0.06)\"),\n    @(12, 3, \"0.91 (0.89 to 0.92)\", \"0.9 (0.89 to 0.92)\"),\n    @(18, 3, \"5.02%\", \"5.23%\"),\n    @(19, 3, \"4.49% (4.24% to 4.74%)\", \"4.76% (4.51% to 5.02%)\"),\n    @(20, 3, \"0.89 (0.84 to 0.94)\", \"0.91 (0.86 to 0.96)\"),\n    @(21, 3, \"-0.54% (-0.79% to -0.28%)\", \"-0.47% (-0.72% to -0.22%)\"),\n    @(22, 3, \"-0.18 (-0.26 to -0.1)\", \"-0.14 (-0.21 to -0.07)\"),\n    @(23, 3, \"-0.16 (-0.22 to -0.1)\", \"-0.16 (-0.21 to -0.1)\"),\n    @(25, 3, \"0.89 (0.88 to 0.9)\", \"0.88 (0.87 to 0.89)\"),\n    @(27, 3, \"0.04 (0.03 to 0.04)\", \"0.04 (0.04 to 0.04)\")\n)\n\nforeach ($edit in $edits) {\n    $rowIndex = $edit[0]\n    $colIndex = $edit[1]\n    $oldText = $edit[2]\n    $newText = $edit[3]\n\n    $cell = $tbl.Cell($rowIndex, $colIndex)\n    $cellRange = $cell.Range\n\n    # Trim the trailing end-of-cell marker characters before comparing.\n    $currentText = $cellRange.Text.TrimEnd([char]7, [char]13)\n\n    if ($currentText -ne $oldText) {\n        throw \"Unexpected text in table cell (row $rowIndex, col $colIndex): got '$currentText', expected '$oldText'\"\n    }\n\n    $cellRange.Text = $newText\n}\n"}
